# Fix data errors for citations and language toggle (NCIOCPL#124)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Correct the Spanish press-release path for the deep-learning cervical
# cancer screening article (it was missing the "-deteccion" suffix).
$ws1.Range("A9").Value = "espanol/noticias/comunicados-de-prensa/2019/aprendizaje-profundo-cancer-cuello-uterino-examenes-de-deteccion"

# Widen column A so the longer path values are readable.
$ws1.Columns.Item(1).ColumnWidth = 99.66666666666667

# Re-fit row 7's height now that column A is wider (its wrapped text no
# longer needs the extra height it had before).
$ws1.Rows.Item(7).AutoFit()

# Move the active selection to A9, matching where the edit was made.
[void]$ws1.Range("A9").Select()
